# Daily attendance processing - 2025-12-26 07:10:09
# Normalize the "Recorded By" (column G) entries: reverse the order of the
# comma-separated recorder names/emails so the most recent recorder is
# listed first, except for entries that include an admin@admin.com entry
# which are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notmatch ",") { continue }
    if ($val -match "admin@admin.com") { continue }

    $parts = $val -split ",\s*"
    $count = $parts.Count
    $revParts = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $revParts += $parts[$i]
    }

    $cell.Value2 = [string]::Join(", ", $revParts)
}
